# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" right after "总计" (before the current
#    "2022-Q2" sheet), duplicated from "2022-Q2" so it inherits identical
#    formatting, then overwrite its values with the new quarter's data
#    (3 fund rows instead of 1).
#  - Insert a new row into "总计" summarizing the new quarter, shifting the
#    older quarters down.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by duplicating the current "2022-Q2"
#    sheet (2nd tab) so all styles/structure match exactly, placed before
#    it, then renamed.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Row 2: 970042 / 国海量化优选一年持有股票C
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "970042"
$q3Sheet.Range("C2").Value = "国海量化优选一年持有股票C"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "7.16"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "87.31"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "0.36"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0258"
$q3Sheet.Range("H2").Value = 2

# Add row 3 (copy format of row 2) -> 004209 / 大成智惠量化多策略灵活配置混合
$q3Sheet.Range("A2:H2").Copy()
$q3Sheet.Range("A3:H3").PasteSpecial($xlPasteFormats)
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").NumberFormat = "@"
$q3Sheet.Range("B3").Value = "004209"
$q3Sheet.Range("C3").Value = "大成智惠量化多策略灵活配置混合"
$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "0.55"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "90.85"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "4.50"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0248"
$q3Sheet.Range("H3").Value = 9

# Add row 4 (copy format of row 2) -> 970041 / 国海量化优选一年持有股票A
$q3Sheet.Range("A2:H2").Copy()
$q3Sheet.Range("A4:H4").PasteSpecial($xlPasteFormats)
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").NumberFormat = "@"
$q3Sheet.Range("B4").Value = "970041"
$q3Sheet.Range("C4").Value = "国海量化优选一年持有股票A"
$q3Sheet.Range("D4").NumberFormat = "@"
$q3Sheet.Range("D4").Value = "0.63"
$q3Sheet.Range("E4").NumberFormat = "@"
$q3Sheet.Range("E4").Value = "87.31"
$q3Sheet.Range("F4").NumberFormat = "@"
$q3Sheet.Range("F4").Value = "0.36"
$q3Sheet.Range("G4").NumberFormat = "@"
$q3Sheet.Range("G4").Value = "0.0023"
$q3Sheet.Range("H4").Value = 2

# ---------------------------------------------------------------------
# 2) Update "总计" sheet: insert a new row 2 for 2022-Q3, shift the rest
#    down, and keep the running index in column A sequential.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()

# The inserted row copies the format of the header row above it (bold),
# not the plain data rows below -- re-copy the correct formatting for
# the whole new row from row 3 (still a plain, unmodified data row).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial($xlPasteFormats)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.05

# Re-sequence column A (0,1,2,...) now that a row was inserted.
for ($i = 0; $i -le 6; $i++) {
    $rowNum = $i + 3
    $totalSheet.Range("A$rowNum").Value = $i + 1
}

$excel.CutCopyMode = $false
